# Update the marksheet's "Marking" and "Total" correct-answer counts,
# and the corresponding correct/total score display.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: number of right answers used for marking row (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Total row: total correct marks (B12): 54 -> 90
$ws.Range("B12").Value = 90

# Total row: correct/total marks display (E12): "52/84" -> "90/140"
$ws.Range("E12").Value = "90/140"
